$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7693683
$ws.Range("J19").Value = 10001483
$ws.Range("L19").Value = 10001483
$ws.Range("N19").Value = -10001833
$ws.Range("H33").Value = 327.66666
$ws.Range("I33").Value = 213.2
$ws.Range("J33").Value = 900
$ws.Range("K33").Value = 213.2
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = 15.80000000000001
$ws.Range("N33").Value = -1358
$ws.Range("H40").Value = 166669150
$ws.Range("J40").Value = 166669150
$ws.Range("L40").Value = 166669150
$ws.Range("N40").Value = -166669500
$ws.Range("H53").Value = 50.5
$ws.Range("I53").Value = 48.8
$ws.Range("K53").Value = 48.8
$ws.Range("M53").Value = 588.2
$ws.Range("H118").Value = 609.0833
$ws.Range("I118").Value = 609.0833
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1827.2499
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -170.2499
$ws.Range("N118").ClearContents()
$ws.Range("H138").Value = 6635.28
$ws.Range("I138").Value = 4241.2
$ws.Range("J138").Value = 6761.284
$ws.Range("K138").Value = 12723.6
$ws.Range("L138").Value = 20283.852
$ws.Range("M138").Value = -7583.599999999999
$ws.Range("N138").Value = -30563.852
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 9008
$ws.Range("I22").Value = 9008
$ws.Range("K22").Value = 9008
$ws.Range("M22").Value = -8709
$ws.Range("H74").Value = 2920.842
$ws.Range("I74").Value = 2907
$ws.Range("K74").Value = 2907
$ws.Range("M74").Value = -2033
$ws.Range("H77").Value = 2920.842
$ws.Range("I77").Value = 2907
$ws.Range("K77").Value = 14535
$ws.Range("M77").Value = -10167
$ws.Range("H97").Value = 1763.2727
$ws.Range("I97").Value = 1370.5238
$ws.Range("J97").Value = 10011
$ws.Range("K97").Value = 1370.5238
$ws.Range("L97").Value = 10011
$ws.Range("M97").Value = -874.5237999999999
$ws.Range("N97").Value = -11003
$ws.Range("H132").Value = 2707605.5
$ws.Range("I132").Value = 4860.8184
$ws.Range("K132").Value = 14582.4552
$ws.Range("M132").Value = -12052.4552
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 4814
$ws.Range("I82").Value = 4814
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 4814
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -4431
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 4814
$ws.Range("I85").Value = 4814
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 4814
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -3488
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 559964.3
$ws.Range("I86").Value = 836547.3
$ws.Range("K86").Value = 836547.3
$ws.Range("M86").Value = -835424.3
$ws.Range("H89").Value = 559964.3
$ws.Range("I89").Value = 836547.3
$ws.Range("K89").Value = 4182736.5
$ws.Range("M89").Value = -4177120.5
$ws.Range("H134").Value = 4349611
$ws.Range("J134").Value = 33335166
$ws.Range("L134").Value = 100005498
$ws.Range("N134").Value = -100010568
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35718964
$ws.Range("I31").Value = 41670396
$ws.Range("K31").Value = 41670396
$ws.Range("M31").Value = -41670101
$ws.Range("H34").Value = 35718964
$ws.Range("I34").Value = 41670396
$ws.Range("K34").Value = 41670396
$ws.Range("M34").Value = -41670194
$ws.Range("H62").Value = 30795.6
$ws.Range("I62").Value = 29999
$ws.Range("J62").Value = 30994.75
$ws.Range("K62").Value = 29999
$ws.Range("L62").Value = 30994.75
$ws.Range("M62").Value = -29375
$ws.Range("N62").Value = -32242.75
$ws.Range("H65").Value = 30795.6
$ws.Range("I65").Value = 29999
$ws.Range("J65").Value = 30994.75
$ws.Range("K65").Value = 149995
$ws.Range("L65").Value = 154973.75
$ws.Range("M65").Value = -146875
$ws.Range("N65").Value = -161213.75
$ws.Range("H74").Value = 75000
$ws.Range("J74").Value = 75000
$ws.Range("L74").Value = 75000
$ws.Range("N74").Value = -76748
$ws.Range("H77").Value = 75000
$ws.Range("J77").Value = 75000
$ws.Range("L77").Value = 225000
$ws.Range("N77").Value = -233736
$ws.Range("H80").Value = 70000
$ws.Range("J80").Value = 70000
$ws.Range("L80").Value = 70000
$ws.Range("N80").Value = -72246
$ws.Range("H83").Value = 70000
$ws.Range("J83").Value = 70000
$ws.Range("L83").Value = 210000
$ws.Range("N83").Value = -221232
$ws.Range("H132").Value = 2355.9524
$ws.Range("J132").Value = 4166.5
$ws.Range("L132").Value = 12499.5
$ws.Range("N132").Value = -17559.5
$ws.Range("H141").Value = 527364.1
$ws.Range("J141").Value = 555394.5
$ws.Range("L141").Value = 555394.5
$ws.Range("N141").Value = -565754.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2221.7778
$ws.Range("I75").Value = 2757.6
$ws.Range("J75").Value = 1552
$ws.Range("K75").Value = 8272.799999999999
$ws.Range("L75").Value = 4656
$ws.Range("M75").Value = -7274.799999999999
$ws.Range("N75").Value = -6652
$ws.Range("H78").Value = 2221.7778
$ws.Range("I78").Value = 2757.6
$ws.Range("J78").Value = 1552
$ws.Range("K78").Value = 24818.4
$ws.Range("L78").Value = 13968
$ws.Range("M78").Value = -19826.4
$ws.Range("N78").Value = -23952
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 913181.75
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25372
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 24000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78864
$ws.Range("H70").Value = 9577.6
$ws.Range("I70").Value = 8999.5
$ws.Range("J70").Value = 9963
$ws.Range("K70").Value = 8999.5
$ws.Range("L70").Value = 9963
$ws.Range("M70").Value = -8729.5
$ws.Range("N70").Value = -10503
$ws.Range("H73").Value = 9577.6
$ws.Range("I73").Value = 8999.5
$ws.Range("J73").Value = 9963
$ws.Range("K73").Value = 8999.5
$ws.Range("L73").Value = 9963
$ws.Range("M73").Value = -8063.5
$ws.Range("N73").Value = -11835
$ws.Range("H97").Value = 1769.6
$ws.Range("I97").Value = 1620.1765
$ws.Range("J97").Value = 2616.3333
$ws.Range("K97").Value = 1620.1765
$ws.Range("L97").Value = 2616.3333
$ws.Range("M97").Value = -1124.1765
$ws.Range("N97").Value = -3608.3333
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6607.75
$ws.Range("I40").Value = 5264.476
$ws.Range("K40").Value = 5264.476
$ws.Range("M40").Value = -5128.476
$ws.Range("H55").Value = 1117.2703
$ws.Range("I55").Value = 657.4286
$ws.Range("J55").Value = 1720.8125
$ws.Range("K55").Value = 657.4286
$ws.Range("L55").Value = 1720.8125
$ws.Range("M55").Value = -484.4286
$ws.Range("N55").Value = -2066.8125
$ws.Range("H132").Value = 3547.1794
$ws.Range("I132").Value = 2466.5557
$ws.Range("K132").Value = 7399.6671
$ws.Range("M132").Value = -4869.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 19999
$ws.Range("J31").Value = 19999
$ws.Range("L31").Value = 19999
$ws.Range("N31").Value = -20695
$ws.Range("H86").Value = 90325
$ws.Range("J86").Value = 90325
$ws.Range("L86").Value = 90325
$ws.Range("N86").Value = -92571
$ws.Range("H89").Value = 90325
$ws.Range("J89").Value = 90325
$ws.Range("L89").Value = 451625
$ws.Range("N89").Value = -462857
$ws.Range("H113").Value = 738
$ws.Range("I113").Value = 807
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 2421
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -251
$ws.Range("N113").Value = -6140
$ws.Range("H126").Value = 5227.744
$ws.Range("I126").Value = 3974.5144
$ws.Range("K126").Value = 11923.5432
$ws.Range("M126").Value = -9453.5432
$ws.Range("H136").Value = 204313.72
$ws.Range("I136").Value = 4229.4893
$ws.Range("J136").Value = 3338966.8
$ws.Range("K136").Value = 12688.4679
$ws.Range("L136").Value = 10016900.4
$ws.Range("M136").Value = -10138.4679
$ws.Range("N136").Value = -10022000.4
